# Updated cryptos list on Fri Mar 31 23:29:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.571.77"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.16"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.48"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5402"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +6.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07768"
$ws.Range("E9").Value = "  +4.62%  "
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.01"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +3.66%  "
$ws.Range("E13").Value = "  +3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.636"
$ws.Range("E14").Value = "  +5.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.001"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.829.16"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001092"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06590"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.574.46"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  +8.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.86"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.14"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.458"
$ws.Range("E28").Value = "  +6.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.040.60"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.46"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.138"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1120"
$ws.Range("E32").Value = "  +5.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.712"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07533"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2258"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02359"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.980"
$ws.Range("E38").Value = "  +6.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.219"
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.43"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6318"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.193"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.409"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.47"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5911"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.710"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.34"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.004"
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("E51").Value = "  +1.43%  "
